$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 131.5
$ws.Range("J5").Value = 174.66667
$ws.Range("L5").Value = 174.66667
$ws.Range("N5").Value = -404.66667
$ws.Range("H11").Value = 59.625
$ws.Range("I11").Value = 59.625
$ws.Range("K11").Value = 59.625
$ws.Range("M11").Value = 80.375
$ws.Range("H33").Value = 220.25
$ws.Range("I33").Value = 201.66667
$ws.Range("K33").Value = 201.66667
$ws.Range("M33").Value = 27.33332999999999
$ws.Range("H70").Value = 1833.1666
$ws.Range("J70").Value = 1833.1666
$ws.Range("L70").Value = 5499.4998
$ws.Range("N70").Value = -6039.4998
$ws.Range("H73").Value = 1833.1666
$ws.Range("J73").Value = 1833.1666
$ws.Range("L73").Value = 5499.4998
$ws.Range("N73").Value = -7371.4998
$ws.Range("H103").Value = 1099.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1099.5
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").Value = 3298.5
$ws.Range("N103").Value = -4470.5
$ws.Range("H106").Value = 2375
$ws.Range("I106").Value = 2375
$ws.Range("K106").Value = 2375
$ws.Range("M106").Value = -1744
$ws.Range("H116").Value = 6416.25
$ws.Range("J116").Value = 6499.5
$ws.Range("L116").Value = 6499.5
$ws.Range("N116").Value = -13383.5
$ws.Range("H137").Value = 1626.95
$ws.Range("I137").Value = 1184.625
$ws.Range("J137").Value = 1921.8334
$ws.Range("K137").Value = 3553.875
$ws.Range("L137").Value = 5765.5002
$ws.Range("M137").Value = -1003.875
$ws.Range("N137").Value = -10865.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1199.125
$ws.Range("I45").Value = 1168.2858
$ws.Range("K45").Value = 1168.2858
$ws.Range("M45").Value = -791.2858000000001
$ws.Range("H61").Value = 2828.3076
$ws.Range("I61").Value = 2828.3076
$ws.Range("K61").Value = 2828.3076
$ws.Range("M61").Value = -2616.3076
$ws.Range("H110").Value = 7401372
$ws.Range("J110").Value = 1500
$ws.Range("L110").Value = 1500
$ws.Range("N110").Value = -5590
$ws.Range("H122").Value = 1364.04
$ws.Range("I122").Value = 1095.5454
$ws.Range("K122").Value = 3286.6362
$ws.Range("M122").Value = -836.6361999999999
$ws.Range("H136").Value = 2828.3076
$ws.Range("I136").Value = 2828.3076
$ws.Range("K136").Value = 8484.9228
$ws.Range("M136").Value = -5934.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1474.5
$ws.Range("I94").Value = 1474.5
$ws.Range("K94").Value = 1474.5
$ws.Range("M94").Value = -1023.5
$ws.Range("H105").Value = 2334.9
$ws.Range("I105").Value = 2356.125
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 2356.125
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = -609.125
$ws.Range("N105").Value = -5744
$ws.Range("H134").Value = 2302.3572
$ws.Range("I134").Value = 2210.2307
$ws.Range("K134").Value = 6630.6921
$ws.Range("M134").Value = -4095.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 86166.664
$ws.Range("I23").Value = 86166.664
$ws.Range("K23").Value = 86166.664
$ws.Range("M23").Value = -85926.664
$ws.Range("H27").Value = 86166.664
$ws.Range("I27").Value = 86166.664
$ws.Range("K27").Value = 86166.664
$ws.Range("M27").Value = -85974.664
$ws.Range("H94").Value = 73137.19
$ws.Range("I94").Value = 140843.5
$ws.Range("K94").Value = 140843.5
$ws.Range("M94").Value = -140392.5
$ws.Range("H107").Value = 503
$ws.Range("I107").Value = 411
$ws.Range("J107").Value = 595
$ws.Range("K107").Value = 411
$ws.Range("L107").Value = 595
$ws.Range("M107").Value = 1509
$ws.Range("N107").Value = -4435
$ws.Range("H132").Value = 4857
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1880.3334
$ws.Range("I38").Value = 2114.375
$ws.Range("K38").Value = 6343.125
$ws.Range("M38").Value = -5996.125
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("H109").Value = 868.93335
$ws.Range("I109").Value = 287.23077
$ws.Range("K109").Value = 861.69231
$ws.Range("M109").Value = 178.30769
$ws.Range("H114").Value = 3270
$ws.Range("I114").Value = 300
$ws.Range("J114").Value = 4012.5
$ws.Range("K114").Value = 900
$ws.Range("L114").Value = 12037.5
$ws.Range("M114").Value = 2354
$ws.Range("N114").Value = -18545.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 90077
$ws.Range("I62").Value = 90077
$ws.Range("K62").Value = 90077
$ws.Range("M62").Value = -89391
$ws.Range("H65").Value = 90077
$ws.Range("I65").Value = 90077
$ws.Range("K65").Value = 270231
$ws.Range("M65").Value = -266799
$ws.Range("H70").Value = 5269.3335
$ws.Range("I70").Value = 5269.3335
$ws.Range("K70").Value = 5269.3335
$ws.Range("M70").Value = -4999.3335
$ws.Range("H73").Value = 5269.3335
$ws.Range("I73").Value = 5269.3335
$ws.Range("K73").Value = 5269.3335
$ws.Range("M73").Value = -4333.3335
$ws.Range("H97").Value = 555
$ws.Range("I97").Value = 555
$ws.Range("K97").Value = 555
$ws.Range("M97").Value = -59
$ws.Range("H102").Value = 2666.4546
$ws.Range("I102").Value = 2333.1
$ws.Range("K102").Value = 2333.1
$ws.Range("M102").Value = -711.0999999999999
$ws.Range("H122").Value = 1794.7142
$ws.Range("I122").Value = 1794.7142
$ws.Range("K122").Value = 5384.142599999999
$ws.Range("M122").Value = -2934.142599999999
$ws.Range("H132").Value = 2949.6
$ws.Range("I132").Value = 2949.6
$ws.Range("K132").Value = 8848.799999999999
$ws.Range("M132").Value = -6318.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 4550
$ws.Range("I18").Value = 4550
$ws.Range("K18").Value = 4550
$ws.Range("M18").Value = -4378
$ws.Range("H99").Value = 9800
$ws.Range("I99").Value = 9800
$ws.Range("K99").Value = 9800
$ws.Range("M99").Value = -6805

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4055.6667
$ws.Range("I126").Value = 3838.8948
$ws.Range("J126").Value = 4570.5
$ws.Range("K126").Value = 11516.6844
$ws.Range("L126").Value = 13711.5
$ws.Range("M126").Value = -9046.6844
$ws.Range("N126").Value = -18651.5
$ws.Range("H132").Value = 3997.2
$ws.Range("I132").Value = 6999.5
$ws.Range("J132").Value = 1995.6666
$ws.Range("K132").Value = 20998.5
$ws.Range("L132").Value = 5986.9998
$ws.Range("M132").Value = -18468.5
$ws.Range("N132").Value = -11046.9998
